$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Kabupaten/Kota sync (feat: sinkronisasi lokasi kabupaten) ---
# M2 "Kabupaten / Kota" value stays "KOTA BOGOR" (unchanged).
# N2 "Kecamatan" changes from "GUNUNG PUTRI" to "BOGOR SELATAN".
# O2 "Desa / Kelurahan" changes from "CIANGSANA" to "CIPAKU".
$ws.Range("N2").Value = "BOGOR SELATAN"
$ws.Range("O2").Value = "CIPAKU"

# --- Row 2: "Tgl. IMB Lama" (C2) and "Tgl. Register" (E2) ---------------
# These were numeric dates (40097) and become literal text "2009-11-10".
# Pre-set the number format to Text ("@") so the literal-looking date
# string is NOT reinterpreted as a date serial number by the engine.

# -4108 = xlCenter, -4160 = xlTop
$ws.Range("C2").HorizontalAlignment = -4108
$ws.Range("C2").VerticalAlignment = -4160
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2009-11-10"

$ws.Range("E2").VerticalAlignment = -4160
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2009-11-10"

# --- Selection / view -----------------------------------------------------
# Active cell moves from M2 to F2, and the frozen/top-left scroll anchor
# resets back to the default (A1) instead of E1.
$ws.Range("A1").Select()
$ws.Range("F2").Select()
